# Updating filtered feeds from workflow
# Adds a new row (40) to the "Filtered Feeds" sheet for a new article:
#   link:     https://www.360dx.com/cancer/roche-nabs-ce-ivdr-marking-her2-cdx-assay-breast-biliary-tract-cancer
#   keywords: CDx
#   title:    Roche Nabs CE-IVDR Marking for HER2 CDx Assay for Breast, Biliary Tract Cancer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink = "https://www.360dx.com/cancer/roche-nabs-ce-ivdr-marking-her2-cdx-assay-breast-biliary-tract-cancer"
$newKeywords = "CDx"
$newTitle = "Roche Nabs CE-IVDR Marking for HER2 CDx Assay for Breast, Biliary Tract Cancer"

$newRow = 40

$ws.Range("A" + $newRow).Value = $newLink
$ws.Range("B" + $newRow).Value = $newKeywords
$ws.Range("C" + $newRow).Value = $newTitle

# Add the hyperlink for the new link cell, matching the style used by the
# other link cells in column A.
$ws.Hyperlinks.Add($ws.Range("A" + $newRow), $newLink)
$ws.Range("A" + $newRow).Style = $ws.Range("A" + ($newRow - 1)).Style

$wb.Save()
